$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 22, pushing the existing rows 22:35 down to 23:36.
$ws.Rows.Item(22).Insert()

# Populate the newly inserted row 22 with the new weekly price entry.
$ws.Range("A22").Value = 7
$ws.Range("B22").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C22").Value = "Ñuble"
$ws.Range("D22").Value = Get-Date -Year 2022 -Month 3 -Day 31 -Hour 0 -Minute 0 -Second 0
$ws.Range("E22").Value = 16
$ws.Range("F22").Value = 100112040
$ws.Range("G22").Value = "Cilantro"
$ws.Range("H22").Value = "Sin especificar"
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 600
$ws.Range("L22").Value = 650
$ws.Range("M22").Value = 625
$ws.Range("N22").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O22").Value = "Provincia de Diguillín"
$ws.Range("P22").Value = 625
$ws.Range("Q22").Value = 1
$ws.Range("R22").Value = "Hortaliza"
